$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -6.221199999999993
$ws.Range("B3").Value = 5.956099999999989
$ws.Range("D5").Value = -8.377899999999993
$ws.Range("E5").Value = 11.86639999999999
$ws.Range("E9").Value = 14.73370000000001
$ws.Range("E11").Value = 13.89
$ws.Range("B14").Value = 8.845600000000005
$ws.Range("B16").Value = 9.429900000000005
$ws.Range("D16").Value = -7.654700000000004
$ws.Range("E17").Value = 13.77650000000001
$ws.Range("B21").Value = 5.741899999999995
$ws.Range("E21").Value = 13.3035
$ws.Range("B23").Value = 5.460200000000002
$ws.Range("B25").Value = 5.820699999999992
